# Mark a few more Filelist rows as translated ("ok") in column B,
# matching the commit "Marked a few strings for translation."
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("B402:B446").Value = "ok"

# Leave the selection where the author ended up after filling the column.
[void]$ws.Range("B447").Select()
